$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.031.86"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "2.624.52"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").Value = "2.622.51"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("E10").Value = "  +10.46%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.355"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("D16").Value = "3.102.52"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "67.865.37"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "2.627.26"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "369.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.70%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").Value = "2.730.96"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "575.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.40%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -3.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.25%  "
$ws.Range("E43").Value = "  -4.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  -6.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.98%  "
